$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.209367632865906
$ws.Range("B1").Value = 2.541895389556885
$ws.Range("C1").Value = 9.365715026855469
$ws.Range("D1").Value = 2.073488473892212
$ws.Range("E1").Value = 1.195730447769165
